# Apply updated dSF (column F) values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -5
    6  = 3
    9  = -1
    11 = 0
    12 = 0
    14 = 3
    17 = 4
    19 = -1
    22 = -1
    26 = 0
    49 = 1
    52 = 0
    56 = 0
    73 = -1
    74 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
